$wb = $excel.ActiveWorkbook

# --- Sheet that currently has tabSelected="1" (infrastructure, physical sheet2.xml) ---
# It loses the tab-selected flag; the "time" sheet (physical sheet5.xml) gains it.
# Physical sheet index 5 corresponds to sheet5.xml, which is the sheet holding the
# gross ton-km derivation method rows edited below. Activating it moves the
# tabSelected flag there and updates the workbook's active-tab pointer.
$wsTime = $wb.Worksheets.Item(5)

# Update the gross ton-km projection method values + number format (percentage, 0.0%).
$wsTime.Range("B2:B6").NumberFormat = "0.0%"

$wsTime.Range("B2").Value = 0.8
$wsTime.Range("B3").Value = 0.7
$wsTime.Range("B4").Value = 0.7
$wsTime.Range("B5").Value = 0.6
$wsTime.Range("B6").Value = 0.5

# Make this sheet the active tab and move the selection to B6, matching the new view state.
$wsTime.Activate()
$wsTime.Range("B6").Select()
